$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Text
    if ($val -eq "Davis") {
        $cell.Value = "T"
    } elseif ($val -eq "Student") {
        $cell.Value = "S"
    }
}
